{"js": "// Merge the two runs \"Factura generada cuando el due\u00f1o realiza un pedido de\n// forma \" + \"virtual\" (same formatting) into a single run so the paragraph\n// text reads \"Factura generada cuando el due\u00f1o realiza un pedido de forma\n// virtual\" as one contiguous <w:t>.\n\nconst targetText =\n  \"Factura generada cuando el due\u00f1o realiza un pedido de forma virtual\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.trim() === targetText);\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph to merge runs in.\");\n}\n\n// Replacing the whole paragraph range with the same text (in one shot)\n// collapses the two separately-formatted-but-identical runs into a single\n// run, exactly like Word does when you delete the run boundary.\nconst wholeRange = target.getRange(\"Whole\");\nwholeRange.insertText(targetText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Merge the two runs \"Factura generada cuando el due\u00f1o realiza un pedido de\n# forma \" + \"virtual\" (both share identical run formatting) into a single\n# run so the paragraph text \"Factura generada cuando el due\u00f1o realiza un\n# pedido de forma virtual\" lives in one contiguous <w:t>.\n\n$wdReplaceOne = 1\n\n$targetText = \"Factura generada cuando el due\u00f1o realiza un pedido de forma virtual\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $targetText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $targetText\n\n# Executing Find/Replace (instead of deleting + re-inserting text) keeps the\n# matched run's existing character formatting (<w:rPr>), which is what\n# collapses the two identically-formatted runs into one clean run.\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceOne) | Out-Null\n"}
